$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.299.48"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.99%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.777.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +4.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.96%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5219"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.29%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3604"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.46%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.55"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07379"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.70%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.091"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.80%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.55"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.15%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.072"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.25%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.777.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.27%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.984"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.53%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.43"
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001045"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06414"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.96%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9999"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.10%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.90%  "

# Row 22
$ws.Range("E22").Value = "  +4.70%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.396.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.21%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.31"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.60%  "

# Row 25
$ws.Range("E25").Value = "  -1.88%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.71%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.96%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.348"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +14.15%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.982.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.35%  "

# Row 30
$ws.Range("E30").Value = "  +2.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.062"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.76%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09781"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.06%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.557"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.27%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.606"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.72%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02236"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.60%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05976"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.00%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "11.20"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.90%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.850"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.84%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2027"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.43%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6150"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.92%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.427"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.65%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.095"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.95%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.150"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.95%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.99%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5767"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.74%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.630"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.13%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.39%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.891"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.46%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.109"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.75%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06718"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.41%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.04%  "
